$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reports")

$ws.Range("C7").Value = "3)"
$ws.Range("D7").Value = "Para la tabla de WorkOrderDetails el registro estara completo Completed = 1 cuado"
$ws.Range("D8").Value = "Quantity = RawMaterial = Machined = Invoiced = Shipped"
$ws.Range("D9").Value = "El TT no intervinene porque ese es solo informativo y sirve para saber si lleva TT o no"

$ws.Range("D10").Select()
